$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing style/format of the data range, then force Text format
# so that numeric-looking price strings (e.g. "1.00", "524.51") are stored as
# text, matching the workbook's existing inline-string convention, not numbers.
$dataRange = $ws.Range("D2:E51")
$origStyle = $ws.Range("D2").Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "57.283.78"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "3.099.61"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "524.51"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "136.63"
$ws.Range("E6").Value = "  -3.47%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.097.62"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.447"
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "0.394"
$ws.Range("E12").Value = "  +2.30%  "
$ws.Range("D13").Value = "3.638.52"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "25.26"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "57.351.66"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "3.099.83"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "5.92"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("D20").Value = "12.35"
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").Value = "344.95"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "67.57"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").Value = "0.499"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "0.0₃0889"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "7.38"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "6.04"
$ws.Range("E32").Value = "  -6.82%  "
$ws.Range("D33").Value = "20.78"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").Value = "4.90"
$ws.Range("E34").Value = "  +6.50%  "
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("D36").Value = "158.27"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "6.05"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").Value = "25.74"
$ws.Range("E38").Value = "  -4.67%  "
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("E40").Value = "  +5.54%  "
$ws.Range("D41").Value = "0.0655"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("D42").Value = "4.09"
$ws.Range("E42").Value = "  +3.39%  "
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "3.140.71"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "2.375.54"
$ws.Range("E45").Value = "  +3.49%  "
$ws.Range("D46").Value = "36.51"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "0.0266"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").Value = "0.972"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("E51").Value = "  -3.96%  "

# Restore the original style/number format so no visible formatting changes remain
$dataRange.Style = $origStyle
